# "Design tests on title bar"
#  1. Bump the cached "Update automatically" date placeholder text
#     (18/08/2022 -> 25/08/2022) on the Slide Master and on every
#     Slide Layout.
#  2. On Slide 1, inside the "Toolbar" group, change the "show items"
#     shape's text from "Show items" (two runs: "Show " @16pt + "items"
#     @10pt) to a single "List" run at 16pt, dropping the now-unused
#     trailing endParaRPr.

$p = $ppt.ActivePresentation

$oldDate = "18/08/2022"
$newDate = "25/08/2022"

# --- 1a. Slide Master date placeholder -------------------------------
$masterShapes = $p.SlideMaster.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $sh = $masterShapes.Item($i)
    if ($sh.Name.StartsWith("Date Placeholder")) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 1b. Every Slide Layout's date placeholder ------------------------
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $cl = $layouts.Item($li)
    $clShapes = $cl.Shapes
    for ($i = 1; $i -le $clShapes.Count; $i++) {
        $sh = $clShapes.Item($i)
        if ($sh.Name.StartsWith("Date Placeholder")) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 2. Slide 1: "Toolbar" group -> "show items" shape ---------------
$slide1 = $p.Slides.Item(1)
$shapes = $slide1.Shapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    if ($sh.Name -eq "Toolbar") {
        $items = $sh.GroupItems
        for ($j = 1; $j -le $items.Count; $j++) {
            $item = $items.Item($j)
            if ($item.Name -eq "show items") {
                $tr = $item.TextFrame.TextRange
                # Replace the whole text; this keeps run #1's formatting
                # (16pt Century Gothic) but the engine always re-adds an
                # endParaRPr after a whole-range assignment, so we redo
                # it as delete-whole-range + insert to drop that
                # trailing endParaRPr (matching the "single run, no
                # endParaRPr" paragraphs used elsewhere in this deck).
                $tr.Text = "List"
                $whole = $tr.Characters(1, 0)
                $whole.Delete()
                $tr.InsertAfter("List") | Out-Null
            }
        }
    }
}
